# Update "想去人数" (F column) figures across sheets to match
# the regenerated data snapshot (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 772
$ws1.Range("F4").Value = 55
$ws1.Range("F7").Value = 14
$ws1.Range("F8").Value = 140
$ws1.Range("F9").Value = 327
$ws1.Range("F10").Value = 442
$ws1.Range("F11").Value = 503
$ws1.Range("F13").Value = 11564
$ws1.Range("F14").Value = 5395

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 100

# --- Sheet "本地生活" ---
# (no data changes in this sheet)

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 772
$ws4.Range("F4").Value = 55
$ws4.Range("F5").Value = 100
$ws4.Range("F9").Value = 14
$ws4.Range("F10").Value = 140
$ws4.Range("F11").Value = 327
$ws4.Range("F12").Value = 442
$ws4.Range("F13").Value = 503
$ws4.Range("F15").Value = 11564
$ws4.Range("F17").Value = 5395
